# The source workbook was a template with a "Folio" header and a
# text-stored placeholder number in A2. This commit "completes" the input
# file: A1 keeps its "Folio" label and A2 receives the real folio number,
# entered as a genuine number (not text) formatted with a plain integer
# number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: header label.
$ws.Range("A1").Value = "Folio"

# A2: the actual folio number, stored as a real number (no thousands
# separator, no decimals) rather than as text.
$ws.Range("A2").Value = 1212300002156
$ws.Range("A2").NumberFormat = "0"

# Match the column sizing / row heights Excel left behind after entering
# the data, and leave the selection on C6 as in the saved file.
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 18.75
$ws.Columns.Item(1).ColumnWidth = 14.6

$ws.Range("C6").Select() | Out-Null
